# Updated cryptos list values (Price / Volume(1h)) per target diff.
# D-column values that look like plain numbers get a leading apostrophe
# so Excel stores them as text (matching the original inline-string cells)
# instead of auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.871.14'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '3.337.02'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''189.21'
$ws.Range("E5").Value = '  +2.33%  '
$ws.Range("D6").Value = '''589.10'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.603'
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").Value = '''6.74'
$ws.Range("E10").Value = '  +2.59%  '
$ws.Range("D11").Value = '''0.414'
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("D12").Value = '3.917.08'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '''28.06'
$ws.Range("E14").Value = '  +2.08%  '
$ws.Range("D15").Value = '68.927.56'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '''0.0000170'
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = '3.311.22'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("D18").Value = '''449.31'
$ws.Range("E18").Value = '  +12.70%  '
$ws.Range("D19").Value = '''5.80'
$ws.Range("D20").Value = '''13.78'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").Value = '''7.82'
$ws.Range("E21").Value = '  +3.01%  '
$ws.Range("D22").Value = '''75.82'
$ws.Range("E22").Value = '  +7.09%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''0.523'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '3.488.93'
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("E26").Value = '  +2.67%  '
$ws.Range("D27").Value = '''0.191'
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("D28").Value = '''9.37'
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +3.06%  '
$ws.Range("D31").Value = '''23.23'
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("D32").Value = '''5.47'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = '''6.93'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +5.92%  '
$ws.Range("D37").Value = '''163.33'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").Value = '''1.93'
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("D39").Value = '''27.12'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").Value = '''4.61'
$ws.Range("E40").Value = '  +2.26%  '
$ws.Range("D41").Value = '''0.796'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").Value = '''6.47'
$ws.Range("E42").Value = '  +2.73%  '
$ws.Range("D43").Value = '2.694.06'
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D45").Value = '''41.16'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D46").Value = '''0.0683'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").Value = '''25.05'
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("D48").Value = '''331.47'
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("D50").Value = '''32.30'
$ws.Range("E50").Value = '  +5.60%  '
$ws.Range("E51").Value = '  +3.38%  '
